# DS_cleanVersion_TEST.xlsx -- "added new version (3.8.1)" edit
#
# Adds two new Xray-integration columns ("xrayTestExecKey", "xrayTestKey")
# in between the existing "testVariantDesc" column and the "dsVal1..3"
# columns, populates the new data row with "n/a", and tidies up the
# header/row formatting that came along with the resave.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room for the two new columns. This shifts the existing
#    dsVal1 / dsVal2 / dsVal3 columns (F,G,H) two places to the right
#    (-> H,I,J), carrying their values/styles/widths with them, which
#    matches the widened <cols> ranges and the H1/I1/J1 shift seen in
#    the diff.
# ---------------------------------------------------------------------
$ws.Columns("F:G").Insert()

# ---------------------------------------------------------------------
# 2. Fill in the headers for the two new columns.
# ---------------------------------------------------------------------
$ws.Range("F1").Value = "xrayTestExecKey"
$ws.Range("G1").Value = "xrayTestKey"

# ---------------------------------------------------------------------
# 3. Fill in the data row under the new headers.
# ---------------------------------------------------------------------
$ws.Range("F2").Value = "n/a"
$ws.Range("G2").Value = "n/a"

# The inserted column carries the column's own default style onto the
# new data cells; the target file leaves F2/G2 with the workbook's
# plain default style (same as the neighbouring D2/E2 cells), so line
# that up explicitly.
$ws.Cells.Item(2, 6).Style = $ws.Cells.Item(2, 4).Style
$ws.Cells.Item(2, 7).Style = $ws.Cells.Item(2, 4).Style

# ---------------------------------------------------------------------
# 4. Row-height tweaks that came along with the resave.
# ---------------------------------------------------------------------
$ws.Rows(1).RowHeight = 28.2
$ws.Rows(2).RowHeight = 15

# ---------------------------------------------------------------------
# 5. Leave the selection where the diff shows it was saved (H2).
# ---------------------------------------------------------------------
$ws.Range("H2").Select()
